# Edit: 'schema iniziale' - restructure DATI FILTRATI section
$d = $word.ActiveDocument

# 1) Find the 'DATI FILTRATI:' heading paragraph (Titolo2) and the following
#    empty Titolo2 paragraph by scanning for the known text, so the script
#    is resilient to any paragraph-count differences.
$titleIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.TrimEnd([char]13, [char]7) -eq "DATI FILTRATI:") {
        $titleIdx = $i
        break
    }
}
if ($titleIdx -eq 0) { throw "Could not locate 'DATI FILTRATI:' paragraph" }

# 2) Replace the heading run with the split-run 'DATI DA FILTRARE:' version.
$titlePara = $d.Paragraphs.Item($titleIdx)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titolo2"/></w:pPr><w:r><w:t xml:space="preserve">DATI </w:t></w:r><w:r><w:t>DA FILTRARE</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.Range.InsertXML($titleXml)

# 3) Replace the following empty 'Titolo2' paragraph with the new block of
#    paragraphs describing the database tables, plus a recreated copy of the
#    original 'DATI FILTRATI:' heading + spacer that now follows the new block.
$emptyIdx = $titleIdx + 1
$emptyPara = $d.Paragraphs.Item($emptyIdx)
$chunkXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Tabella </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>countries</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Colonne</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>: id, country;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Tabella</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>program_languages</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Colonne</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>: id, language;</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Tabella</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>trash_tweet</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>Colonne</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">: id, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>id_country</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve">, date, text, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>id_lang</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
        <w:t xml:space="preserve"> retweet, </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-GB"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Titolo2"/>
      </w:pPr>
      <w:r>
        <w:t>DATI FILTRATI:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Titolo2"/>
      </w:pPr>
    </w:p>
    </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$emptyPara.Range.InsertXML($chunkXml)

# 4) The original '_GoBack' bookmark that sat in the 'Colonne dal database...'
#    paragraph has now been recreated further up (inside the new trash_tweet
#    column list); remove the stale original so only one '_GoBack' remains.
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBm.Delete()
} catch {
}

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
